$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at the top; existing data (old rows 1-7) shifts down to rows 5-11
$ws.Range("A1:A4").EntireRow.Insert()

# Row 1: new header labels
$ws.Range("A1").Value = "Contrast"
$ws.Range("B1").Value = "Correlation"
$ws.Range("C1").Value = "Energy"
$ws.Range("D1").Value = "Homogeneity"
$ws.Range("E1").Value = "Actual Condition"

# Row 2: new data row
$ws.Range("A2").Value = 0.004827386532050235
$ws.Range("B2").Value = 0.9837085928765761
$ws.Range("C2").Value = 0.6988810257581656
$ws.Range("D2").Value = 0.9975863067339747
$ws.Range("E2").Value = "N"

# Row 3: new data row
$ws.Range("A3").Value = 0.010702044568282389
$ws.Range("B3").Value = 0.9441616319776373
$ws.Range("C3").Value = 0.7977513617728199
$ws.Range("D3").Value = 0.9946489777158589
$ws.Range("E3").Value = "B"

# Row 4: new data row
$ws.Range("A4").Value = 0.005918612862137599
$ws.Range("B4").Value = 0.9825819331507765
$ws.Range("C4").Value = 0.6543190725903568
$ws.Range("D4").Value = 0.9970406935689312
$ws.Range("E4").Value = "B"

# Update the active selection to match the committed workbook state
$ws.Range("E18").Select() | Out-Null
